$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1736.25
$ws.Range("I15").Value = 1736.25
$ws.Range("K15").Value = 5208.75
$ws.Range("M15").Value = -5039.75
$ws.Range("H19").Value = 1327
$ws.Range("I19").Value = 1290.9231
$ws.Range("K19").Value = 1290.9231
$ws.Range("M19").Value = -1115.9231
$ws.Range("H53").Value = 281.75
$ws.Range("I53").Value = 248.54546
$ws.Range("K53").Value = 248.54546
$ws.Range("M53").Value = 388.45454
$ws.Range("H98").Value = 798.8
$ws.Range("I98").Value = 798.8
$ws.Range("K98").Value = 798.8
$ws.Range("M98").Value = 699.2
$ws.Range("H122").Value = 798.8
$ws.Range("I122").Value = 798.8
$ws.Range("K122").Value = 2396.4
$ws.Range("M122").Value = 53.60000000000036
$ws.Range("H132").Value = 1518.3334
$ws.Range("I132").Value = 1422.2
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 4266.6
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -1736.6
$ws.Range("N132").Value = -11057
$ws.Range("H138").Value = 4038.9143
$ws.Range("I138").Value = 2743.9167
$ws.Range("J138").Value = 4714.5654
$ws.Range("K138").Value = 8231.750100000001
$ws.Range("L138").Value = 14143.6962
$ws.Range("M138").Value = -3091.750100000001
$ws.Range("N138").Value = -24423.6962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 150.33333
$ws.Range("I5").Value = 125.5
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 125.5
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -13.5
$ws.Range("N5").Value = -424
$ws.Range("H32").Value = 7141.0527
$ws.Range("I32").Value = 5922.4116
$ws.Range("K32").Value = 5922.4116
$ws.Range("M32").Value = -5635.4116
$ws.Range("H61").Value = 1754.5
$ws.Range("I61").Value = 1680.1
$ws.Range("K61").Value = 1680.1
$ws.Range("M61").Value = -1468.1
$ws.Range("H74").Value = 457
$ws.Range("J74").Value = 457
$ws.Range("L74").Value = 457
$ws.Range("N74").Value = -2205
$ws.Range("H77").Value = 457
$ws.Range("J77").Value = 457
$ws.Range("L77").Value = 2285
$ws.Range("N77").Value = -11021
$ws.Range("H97").Value = 314.3
$ws.Range("I97").Value = 314.3
$ws.Range("K97").Value = 314.3
$ws.Range("M97").Value = 181.7
$ws.Range("H136").Value = 1754.5
$ws.Range("I136").Value = 1680.1
$ws.Range("K136").Value = 5040.299999999999
$ws.Range("M136").Value = -2490.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150.33333
$ws.Range("I4").Value = 125.5
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 125.5
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -10.5
$ws.Range("N4").Value = -430
$ws.Range("H94").Value = 6767.1055
$ws.Range("I94").Value = 8246.733
$ws.Range("K94").Value = 8246.733
$ws.Range("M94").Value = -7795.733
$ws.Range("H97").Value = 8875.5
$ws.Range("I97").Value = 8875.5
$ws.Range("K97").Value = 8875.5
$ws.Range("M97").Value = -7884.5
$ws.Range("H99").Value = 814.1111
$ws.Range("I99").Value = 827.25
$ws.Range("J99").Value = 709
$ws.Range("K99").Value = 827.25
$ws.Range("L99").Value = 709
$ws.Range("M99").Value = 670.75
$ws.Range("N99").Value = -3705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H22").Value = 2040
$ws.Range("I22").Value = 2040
$ws.Range("K22").Value = 2040
$ws.Range("M22").Value = -1690
$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630
$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184
$ws.Range("H97").Value = 30197
$ws.Range("J97").Value = 30197
$ws.Range("L97").Value = 30197
$ws.Range("N97").Value = -32179
$ws.Range("H107").Value = 728.4167
$ws.Range("I107").Value = 389.7
$ws.Range("K107").Value = 389.7
$ws.Range("M107").Value = 1530.3
$ws.Range("H109").Value = 59994.5
$ws.Range("J109").Value = 59994.5
$ws.Range("L109").Value = 59994.5
$ws.Range("N109").Value = -62074.5
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H132").Value = 4475
$ws.Range("I132").Value = 3746.7144
$ws.Range("K132").Value = 11240.1432
$ws.Range("M132").Value = -8710.143199999999
$ws.Range("H133").Value = 124497.5
$ws.Range("J133").Value = 124497.5
$ws.Range("L133").Value = 124497.5
$ws.Range("N133").Value = -129557.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 933
$ws.Range("I10").Value = 299
$ws.Range("J10").Value = 1250
$ws.Range("K10").Value = 897
$ws.Range("L10").Value = 3750
$ws.Range("M10").Value = -758
$ws.Range("N10").Value = -4028
$ws.Range("H33").Value = 84
$ws.Range("I33").Value = 70.333336
$ws.Range("J33").Value = 125
$ws.Range("K33").Value = 422.000016
$ws.Range("L33").Value = 750
$ws.Range("M33").Value = -139.000016
$ws.Range("N33").Value = -1316
$ws.Range("H109").Value = 875.4
$ws.Range("I109").Value = 844.75
$ws.Range("K109").Value = 2534.25
$ws.Range("M109").Value = -1494.25
$ws.Range("H113").Value = 739.1818
$ws.Range("I113").Value = 499.16666
$ws.Range("K113").Value = 1497.49998
$ws.Range("M113").Value = 672.5000199999999
$ws.Range("H137").Value = 3750
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1777.3334
$ws.Range("I102").Value = 1777.3334
$ws.Range("K102").Value = 1777.3334
$ws.Range("M102").Value = -155.3334
$ws.Range("H107").Value = 1399.25
$ws.Range("I107").Value = 1399.25
$ws.Range("K107").Value = 1399.25
$ws.Range("M107").Value = 520.75
$ws.Range("H113").Value = 1394
$ws.Range("I113").Value = 1394
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1394
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 776
$ws.Range("H122").Value = 1000.875
$ws.Range("I122").Value = 1001.2857
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 3003.8571
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -553.8571000000002
$ws.Range("N122").Value = -7894
$ws.Range("H126").Value = 2061.2222
$ws.Range("I126").Value = 1517.3334
$ws.Range("K126").Value = 4552.0002
$ws.Range("M126").Value = -2082.0002
$ws.Range("H132").Value = 3542.5625
$ws.Range("I132").Value = 3191.6428
$ws.Range("K132").Value = 9574.928400000001
$ws.Range("M132").Value = -7044.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3497.6667
$ws.Range("I7").Value = 2997.2
$ws.Range("K7").Value = 2997.2
$ws.Range("M7").Value = -2885.2
$ws.Range("H46").Value = 2908.25
$ws.Range("I46").Value = 2373
$ws.Range("J46").Value = 3443.5
$ws.Range("K46").Value = 2373
$ws.Range("L46").Value = 3443.5
$ws.Range("M46").Value = -2185
$ws.Range("N46").Value = -3819.5
$ws.Range("H61").Value = 4982.6665
$ws.Range("I61").Value = 4974.5
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 4974.5
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -4772.5
$ws.Range("N61").Value = -5403
$ws.Range("H113").Value = 4982.6665
$ws.Range("I113").Value = 4974.5
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 4974.5
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = -2804.5
$ws.Range("N113").Value = -9339
$ws.Range("H126").Value = 3497.6667
$ws.Range("I126").Value = 2997.2
$ws.Range("K126").Value = 8991.599999999999
$ws.Range("M126").Value = -6521.599999999999
$ws.Range("H132").Value = 4798.643
$ws.Range("J132").Value = 4847.25
$ws.Range("L132").Value = 14541.75
$ws.Range("N132").Value = -19601.75
$ws.Range("H136").Value = 26109.842
$ws.Range("I136").Value = 8735.909
$ws.Range("K136").Value = 26207.727
$ws.Range("M136").Value = -23657.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 75000
$ws.Range("J46").Value = 75000
$ws.Range("L46").Value = 75000
$ws.Range("N46").Value = -75462
$ws.Range("H100").Value = 1073.5
$ws.Range("I100").Value = 1399.5
$ws.Range("K100").Value = 2799
$ws.Range("M100").Value = -2258
$ws.Range("H113").Value = 774.1667
$ws.Range("I113").Value = 795.4
$ws.Range("J113").Value = 668
$ws.Range("K113").Value = 2386.2
$ws.Range("L113").Value = 2004
$ws.Range("M113").Value = -216.1999999999998
$ws.Range("N113").Value = -6344
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 1778.3
$ws.Range("I126").Value = 1809.4615
$ws.Range("J126").Value = 1720.4286
$ws.Range("K126").Value = 5428.3845
$ws.Range("L126").Value = 5161.2858
$ws.Range("M126").Value = -2958.3845
$ws.Range("N126").Value = -10101.2858
$ws.Range("H134").Value = 75000
$ws.Range("J134").Value = 75000
$ws.Range("L134").Value = 225000
$ws.Range("N134").Value = -230070

